$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-09-16"

# Update the row label for September to reflect the new "through" date
$ws.Range("A10").Value = "September (through 09-16)"

# Update September row (row 10) values for years 2015-2021 (columns B-H)
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 29
$ws.Range("D10").Value = 40
$ws.Range("E10").Value = 30
$ws.Range("F10").Value = 35
$ws.Range("G10").Value = 61
$ws.Range("H10").Value = 80

# Update Total row (row 11) values for years 2015-2021 (columns B-H)
$ws.Range("B11").Value = 211
$ws.Range("C11").Value = 410
$ws.Range("D11").Value = 591
$ws.Range("E11").Value = 520
$ws.Range("F11").Value = 384
$ws.Range("G11").Value = 845
$ws.Range("H11").Value = 1150
